# Puff Smith translations.xlsx - "Common: A lot of small improvements; added glow and active flag"
# Adds 8 new "cs" translation rows (glow label/tooltip/table header, mixture
# activate/deactivate buttons + success toasts, mixture preview "active" flag)
# to the bottom of the "Import" sheet's data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New rows to append right after the current last row (626).
# Each entry: row number, translation key (column B), translation text
# (column C), and whether the value (C) was originally typed in before the
# key (B) - matches the author's original shared-string insertion order.
$newRows = @(
    @{ Row = 627; Key = "lab.build.glow.label";              Text = "Žhavení"; ValueFirst = $true },
    @{ Row = 628; Key = "lab.build.glow.label.tooltip";       Text = 'Tato hodnota udává rychlost žhavení; čím vyšší číslo, tím rychleji se spirálka rozžhaví; smyslem je poskytnout náhled, jak moc "divoký" build je.'; ValueFirst = $false },
    @{ Row = 629; Key = "lab.build.table.glow";                Text = "Žhavení"; ValueFirst = $true },
    @{ Row = 630; Key = "lab.mixture.button.activate";         Text = "Aktivovat mix"; ValueFirst = $false },
    @{ Row = 631; Key = "lab.mixture.button.deactivate";       Text = "Deaktivovat mix"; ValueFirst = $false },
    @{ Row = 632; Key = "lab.mixture.deactivated.success";     Text = "Mix [{{data.name}}] byl úspěšně deaktivován."; ValueFirst = $false },
    @{ Row = 633; Key = "lab.mixture.activated.success";       Text = "Mix [{{data.name}}] byl úspěšně aktivován."; ValueFirst = $false },
    @{ Row = 634; Key = "lab.mixture.preview.active";          Text = "Aktivní"; ValueFirst = $false }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $rowRange = $ws.Range("A" + $r + ":C" + $r)
    # Match the formatting already used by every other data row (wrapped
    # text, 10pt font - the "import" cell style).
    $rowRange.WrapText = $true
    $rowRange.Font.Size = 10

    $ws.Range("A" + $r).Value = "cs"
    if ($item.ValueFirst) {
        $ws.Range("C" + $r).Value = $item.Text
        $ws.Range("B" + $r).Value = $item.Key
    } else {
        $ws.Range("B" + $r).Value = $item.Key
        $ws.Range("C" + $r).Value = $item.Text
    }
}

# The tooltip text is long enough to wrap onto two lines - give that row the
# taller row height Excel would compute for it.
$ws.Rows.Item(628).RowHeight = 26.25

# Leave the view the way the author left it: scrolled down to the new rows,
# with the last-edited cell (B629) selected.
$excel.ActiveWindow.ScrollRow = 621
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B629").Select()
